$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Phase 0: stash the two formats already used on the sheet (the plain
# wrap/vcenter style used throughout column B/D, and the quotePrefix
# style already sitting - unused - on C20/C21) onto far-away scratch
# cells. We do this before touching any values so the source cells'
# own formatting is still intact, and so re-applying these formats
# later never causes Excel to mint brand-new style entries.
# ---------------------------------------------------------------------
$ws.Cells.Item(3, 2).Copy()                 # B3 -> wrap/vcenter style
$ws.Cells.Item(1, 26).PasteSpecial(-4122)   # Z1 scratch
$ws.Cells.Item(20, 3).Copy()                # C20 -> quotePrefix style
$ws.Cells.Item(2, 26).PasteSpecial(-4122)   # Z2 scratch
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Phase 1: enter all the new text, in the same left-to-right / top-to-
# bottom order it was originally typed in (step names first, then the
# corresponding values), so new shared-string entries land in the
# expected order.
# ---------------------------------------------------------------------
$ws.Cells.Item(15, 2).Value = "Spreadsheet SpreadsheetResult calc()"

$ws.Cells.Item(16, 2).Value = "Step Name"
$ws.Cells.Item(16, 3).Value = "Value"

$ws.Cells.Item(17, 2).Value = "Step1"
$ws.Cells.Item(18, 2).Value = "Step2"
$ws.Cells.Item(19, 2).Value = "Step3"
$ws.Cells.Item(20, 2).Value = "Step4"
$ws.Cells.Item(21, 2).Value = "Step5"

# Leading apostrophe forces literal text (these all start with "="),
# which lines up with the quotePrefix style applied below.
$ws.Cells.Item(17, 3).Value = '''= for (int java=0;java<10;java++) {} java.lang.Boolean.TRUE;'
$ws.Cells.Item(18, 3).Value = '''= for (int i=0;i<10;i++) { String java = "hello";} java.lang.Boolean.TRUE;'
$ws.Cells.Item(19, 3).Value = '''= if (true) {String java = "Hello";} java.lang.Boolean.TRUE;'
$ws.Cells.Item(20, 3).Value = '''= while (false) {String java = "Hello";} java.lang.Boolean.TRUE;'
$ws.Cells.Item(21, 3).Value = '''= {String java = "Hello";} java.lang.Boolean.TRUE;'

$ws.Cells.Item(22, 2).Value = "Step6"
$ws.Cells.Item(22, 3).Value = "while ("

# ---------------------------------------------------------------------
# Phase 2: re-apply formatting. Column B step names (rows 18-22) and
# the new header row (15) use the normal wrap/vcenter style; column C
# answers (rows 17-22) use the quotePrefix style. Row 16/17's B cells
# and row 16's C cell intentionally keep the default (unstyled) look.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 26).Copy()
$ws.Range("B18:B22").PasteSpecial(-4122)

$ws.Cells.Item(2, 26).Copy()
$ws.Range("C17:C22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Phase 3: clean up the scratch cells.
# ---------------------------------------------------------------------
$ws.Range("Z1:Z2").Clear()

# ---------------------------------------------------------------------
# Phase 4: restore the saved selection.
# ---------------------------------------------------------------------
$ws.Range("C22").Select()
